$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 115; everything currently at 115..190
# shifts down to 116..191 (dimension grows from R190 to R191).
$ws.Rows("115:115").Insert()

# Populate the newly-inserted row 115 with its data.
$ws.Range("A115").Value = 4
$ws.Range("B115").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C115").Value = "Los Lagos"
$ws.Range("D115").Value = 44596
$ws.Range("E115").Value = 10
$ws.Range("F115").Value = 100112032
$ws.Range("G115").Value = "Zapallo italiano"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 200
$ws.Range("K115").Value = 16000
$ws.Range("L115").Value = 16000
$ws.Range("M115").Value = 16000
$ws.Range("N115").Value = "`$/caja 50 unidades"
$ws.Range("O115").Value = "Región de O'Higgins"
$ws.Range("P115").Value = 320
$ws.Range("Q115").Value = 50
$ws.Range("R115").Value = "Hortaliza"
